$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update data values ---

# Row 2 - Venous pH
$ws.Range("B2").Value = 7.38
$ws.Range("C2").Value = 7.32
$ws.Range("D2").Value = 7.4
$ws.Range("E2").Value = 7.54
$ws.Range("F2").Value = 7.54

# Row 3 - Venous pCO2
$ws.Range("B3").Value = 48
$ws.Range("C3").Value = 73
$ws.Range("D3").Value = 74
$ws.Range("E3").Value = 54
$ws.Range("F3").Value = 53

# Row 4 - Arterial pH
$ws.Range("B4").Value = 7.43
$ws.Range("C4").Value = 7.42
$ws.Range("D4").Value = 7.47
$ws.Range("E4").Value = 7.56
$ws.Range("F4").Value = 7.57

# Row 5 - Arterial pCO2
$ws.Range("B5").Value = 43
$ws.Range("C5").Value = 59
$ws.Range("D5").Value = 63
$ws.Range("E5").Value = 51
$ws.Range("F5").Value = 50

# Row 7 - Ventilation
$ws.Range("B7").Value = 6.5
$ws.Range("C7").Value = 12.8
$ws.Range("D7").Value = 9.9
$ws.Range("E7").Value = 4.6

# --- Update row heights (rows 1, 2, 4, 6, 7 shrink from 30.75 to 15.75) ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15.75

# --- Update view: scroll back to top and select F3 instead of A1:F7 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F3").Select()
